$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "C004"
$ws.Range("B5").Value = "J"
$ws.Range("C5").Value = 122356789
$ws.Range("D5").Value = 1234567893713
